$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Activate()

# Insert a new row at 97 (existing rows 97-253 shift down to 98-254) and
# populate it with the new "LFHVM" acronym entry under the "elec" folder.
$ws.Rows.Item(97).Insert()

$ws.Cells.Item(97, 1).Value = "elec"
$ws.Cells.Item(97, 2).Value = "LFHVM"
$ws.Cells.Item(97, 3).Value = "Load Factor Hourly Variance Multiplier"
$ws.Cells.Item(97, 6).Value = "high"

# Match the "high" formatting used elsewhere in column F (e.g. F2) so the
# new row's fill/style lines up with the rest of the "high" rows.
$ws.Range("F2").Copy()
$ws.Range("F97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection to the newly inserted row, matching the edited file.
$ws.Range("A97").Select()
